# Update "想去人数" (F column) figures on the "展览", "演出" and "全部类型"
# sheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) -------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1288
$ws1.Range("F4").Value  = 932
$ws1.Range("F5").Value  = 972
$ws1.Range("F6").Value  = 1735
$ws1.Range("F7").Value  = 382
$ws1.Range("F8").Value  = 1155
$ws1.Range("F9").Value  = 49
$ws1.Range("F10").Value = 6
$ws1.Range("F11").Value = 114
$ws1.Range("F12").Value = 264
$ws1.Range("F13").Value = 44
$ws1.Range("F15").Value = 645
$ws1.Range("F16").Value = 135
$ws1.Range("F17").Value = 90
$ws1.Range("F21").Value = 109
$ws1.Range("F22").Value = 647
$ws1.Range("F27").Value = 844
$ws1.Range("F28").Value = 299
$ws1.Range("F29").Value = 121
$ws1.Range("F30").Value = 27
$ws1.Range("F31").Value = 252
$ws1.Range("F34").Value = 398

# --- Sheet "演出" (sheet2) --------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 311

# --- Sheet "全部类型" (sheet4) ----------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1288
$ws4.Range("F5").Value  = 932
$ws4.Range("F6").Value  = 972
$ws4.Range("F7").Value  = 1735
$ws4.Range("F8").Value  = 382
$ws4.Range("F9").Value  = 1155
$ws4.Range("F10").Value = 49
$ws4.Range("F12").Value = 6
$ws4.Range("F13").Value = 114
$ws4.Range("F14").Value = 264
$ws4.Range("F15").Value = 44
$ws4.Range("F17").Value = 645
$ws4.Range("F18").Value = 135
$ws4.Range("F19").Value = 90
$ws4.Range("F22").Value = 311
$ws4.Range("F29").Value = 109
$ws4.Range("F30").Value = 647
$ws4.Range("F35").Value = 844
$ws4.Range("F36").Value = 299
$ws4.Range("F39").Value = 121
$ws4.Range("F40").Value = 27
$ws4.Range("F41").Value = 252
$ws4.Range("F48").Value = 398
